$d = $word.ActiveDocument

function Set-ParagraphText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $p.Range.Text = $newText
}

# Paragraph 2: intro paragraph about Google Analytics / Zap / Piwik tools
Set-ParagraphText 2 "Strumenti come Google Analytics, Zap o Piwik, tutti permettono di vedere informazioni approfondite sul tuo sito web. Calcola il volume del traffico di ricerca e il traffico potenziale dai mercati in cui al momento non stai mirando. Questi strumenti ti aiuteranno a capire la quantità specifica di traffico locale, a conoscere le fonti del tuo traffico e le informazioni cruciali per un'ottimizzazione SEO, parole chiave e frasi che generano più traffico."

# Paragraph 3: bold heading "In-depth keyword research in target markets"
Set-ParagraphText 3 "Ricerca approfondita di parole chiave nei mercati di destinazione"

# Paragraph 4: keyword research data paragraph
Set-ParagraphText 4 "Con i dati raccolti in precedenza o lavorando con il team di marketing, esegui una ricerca sulle parole chiave che guidano maggiormente il traffico. Cerca specifici modelli di ricerca e tendenze. Puoi utilizzare il pianificatore di parole chiave di Google per scoprire preziose parole chiave utilizzate di frequente."

# Paragraph 7: bold heading "Calculate ROI - Return of Investment"
Set-ParagraphText 7 "Calcola il ROI - Return of Investment"

# Paragraph 8: localization costs paragraph
Set-ParagraphText 8 "I costi di localizzazione si basano sul volume di contenuti che è necessario tradurre e localizzare. È necessario valutare attentamente quali contenuti devono essere realmente tradotti in quali mercati target."

# Paragraph 9: performance indicators paragraph (keeps trailing space)
Set-ParagraphText 9 "Esistono molti indicatori e parametri di rendimento che è possibile utilizzare per misurare i costi e i rendimenti di un sito web localizzato rispetto a un sito Web non localizzato. "

# Paragraph 12: localization decisions paragraph
Set-ParagraphText 12 "Le decisioni relative alla localizzazione devono prendere in considerazione molteplici fattori, tra cui il potenziale dei mercati target e la concorrenza esistente in quei locali. È estremamente importante dotarsi dei numeri e delle metriche corretti che supportano la strategia di localizzazione. Stimare costi e rendimenti con chiari indicatori di prestazione è molto utile per valutare il ROI complessivo di localizzazione."

Write-Output "Done"
